$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 240-241; everything currently at row 240 and below
# shifts down by two rows (so old row 240 becomes row 242, ..., old row 366
# becomes row 368).
$ws.Rows("240:241").Insert()

# --- Row 240 (new weekly record - "Primera") ---
$ws.Range("A240").Value = 3
$ws.Range("B240").Value = "Femacal de La Calera"
$ws.Range("C240").Value = "Coquimbo"
$ws.Range("D240").Value = 44452
$ws.Range("E240").Value = 5
$ws.Range("F240").Value = 100114014
$ws.Range("G240").Value = "Betarraga"
$ws.Range("H240").Value = "Sin especificar"
$ws.Range("I240").Value = "Primera"
$ws.Range("J240").Value = 3100
$ws.Range("K240").Value = 550
$ws.Range("L240").Value = 600
$ws.Range("M240").Value = 576
$ws.Range("N240").Value = "$/paquete 4 unidades"
$ws.Range("O240").Value = "Provincia de Quillota"
$ws.Range("P240").Value = 144
$ws.Range("Q240").Value = 4
$ws.Range("R240").Value = "Hortaliza"

# --- Row 241 (new weekly record - "Segunda") ---
$ws.Range("A241").Value = 3
$ws.Range("B241").Value = "Femacal de La Calera"
$ws.Range("C241").Value = "Coquimbo"
$ws.Range("D241").Value = 44452
$ws.Range("E241").Value = 5
$ws.Range("F241").Value = 100114014
$ws.Range("G241").Value = "Betarraga"
$ws.Range("H241").Value = "Sin especificar"
$ws.Range("I241").Value = "Segunda"
$ws.Range("J241").Value = 1500
$ws.Range("K241").Value = 450
$ws.Range("L241").Value = 450
$ws.Range("M241").Value = 450
$ws.Range("N241").Value = "$/paquete 4 unidades"
$ws.Range("O241").Value = "Provincia de Quillota"
$ws.Range("P241").Value = 112
$ws.Range("Q241").Value = 4
$ws.Range("R241").Value = "Hortaliza"
